# Auto-generated edit script: update cryptos D/E columns per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    # Force text storage so numeric/date-like strings (e.g. '227.41',
    # '1.00', thousands-dotted prices) are NOT coerced to numbers by COM,
    # matching the workbook's original inlineStr cell type. Then restore
    # the default 'Normal' style so no stray number-format style sticks.
    $Range.NumberFormat = '@'
    $Range.Value = $Value
    $Range.Style = 'Normal'
}

# Subscript-three char (U+2083) used in the ShibaInu price; built via
# [char] + interpolation because literal '+' concatenation of a [char]
# with numeric-looking strings performs numeric addition, not text concat.
$sub3 = [string][char]8323

Set-TextValue $ws.Range('D2') '37.845.99'
Set-TextValue $ws.Range('E2') '  -0.79%  '
Set-TextValue $ws.Range('D3') '2.036.19'
Set-TextValue $ws.Range('E3') '  -1.27%  '
Set-TextValue $ws.Range('E4') '  -0.03%  '
Set-TextValue $ws.Range('D5') '227.41'
Set-TextValue $ws.Range('E5') '  -1.20%  '
Set-TextValue $ws.Range('D6') '0.614'
Set-TextValue $ws.Range('E6') '  -0.52%  '
Set-TextValue $ws.Range('D7') '60.15'
Set-TextValue $ws.Range('E7') '  +3.37%  '
Set-TextValue $ws.Range('E8') '  +0.02%  '
Set-TextValue $ws.Range('E9') '  -0.13%  '
Set-TextValue $ws.Range('D10') '0.0818'
Set-TextValue $ws.Range('E10') '  +1.13%  '
Set-TextValue $ws.Range('E11') '  +0.68%  '
Set-TextValue $ws.Range('D12') '14.65'
Set-TextValue $ws.Range('E12') '  +0.20%  '
Set-TextValue $ws.Range('D13') '2.337.13'
Set-TextValue $ws.Range('E13') '  -1.23%  '
Set-TextValue $ws.Range('D14') '21.05'
Set-TextValue $ws.Range('E14') '  +1.53%  '
Set-TextValue $ws.Range('D15') '0.760'
Set-TextValue $ws.Range('E15') '  +0.60%  '
Set-TextValue $ws.Range('D16') '5.23'
Set-TextValue $ws.Range('E16') '  -1.20%  '
Set-TextValue $ws.Range('D17') '2.036.00'
Set-TextValue $ws.Range('E17') '  -1.70%  '
Set-TextValue $ws.Range('D18') '37.797.38'
Set-TextValue $ws.Range('E18') '  -0.63%  '
Set-TextValue $ws.Range('D19') '6.07'
Set-TextValue $ws.Range('E19') '  -2.16%  '
Set-TextValue $ws.Range('D20') '69.83'
Set-TextValue $ws.Range('E20') '  -0.02%  '
Set-TextValue $ws.Range('D21') "0.0${sub3}0824"
Set-TextValue $ws.Range('E21') '  -0.90%  '
Set-TextValue $ws.Range('D22') '225.57'
Set-TextValue $ws.Range('E22') '  +0.35%  '
Set-TextValue $ws.Range('E23') '  +0.04%  '
Set-TextValue $ws.Range('E24') '  -2.20%  '
Set-TextValue $ws.Range('D25') '2.21'
Set-TextValue $ws.Range('E25') '  -2.37%  '
Set-TextValue $ws.Range('D26') '9.27'
Set-TextValue $ws.Range('E26') '  -0.12%  '
Set-TextValue $ws.Range('D27') '165.10'
Set-TextValue $ws.Range('E27') '  -0.32%  '
Set-TextValue $ws.Range('D28') '0.129'
Set-TextValue $ws.Range('E28') '  -4.01%  '
Set-TextValue $ws.Range('D29') '18.92'
Set-TextValue $ws.Range('E29') '  -0.82%  '
Set-TextValue $ws.Range('D30') '1.29'
Set-TextValue $ws.Range('E30') '  -6.83%  '
Set-TextValue $ws.Range('E31') '  +1.47%  '
Set-TextValue $ws.Range('E32') '  -2.67%  '
Set-TextValue $ws.Range('D33') '2.05'
Set-TextValue $ws.Range('E33') '  +3.90%  '
Set-TextValue $ws.Range('E34') '  -2.09%  '
Set-TextValue $ws.Range('D35') '4.48'
Set-TextValue $ws.Range('E35') '  -2.85%  '
Set-TextValue $ws.Range('D36') '6.42'
Set-TextValue $ws.Range('E36') '  +6.66%  '
Set-TextValue $ws.Range('D37') '2.25'
Set-TextValue $ws.Range('E37') '  -5.39%  '
Set-TextValue $ws.Range('E38') '  -1.88%  '
Set-TextValue $ws.Range('E39') '  +0.03%  '
Set-TextValue $ws.Range('D40') '1.539.96'
Set-TextValue $ws.Range('E40') '  +4.11%  '
Set-TextValue $ws.Range('D41') '16.96'
Set-TextValue $ws.Range('E41') '  +0.82%  '
Set-TextValue $ws.Range('E42') '  -0.79%  '
Set-TextValue $ws.Range('D43') '96.90'
Set-TextValue $ws.Range('E43') '  -1.59%  '
Set-TextValue $ws.Range('E44') '  -1.58%  '
Set-TextValue $ws.Range('D45') '0.0921'
Set-TextValue $ws.Range('E45') '  -2.57%  '
Set-TextValue $ws.Range('E46') '  -1.55%  '
Set-TextValue $ws.Range('E47') '  -4.84%  '
Set-TextValue $ws.Range('D48') '1.01'
Set-TextValue $ws.Range('E48') '  -1.78%  '
Set-TextValue $ws.Range('E49') '  -0.29%  '
Set-TextValue $ws.Range('D50') '7.14'
Set-TextValue $ws.Range('E50') '  +0.43%  '
Set-TextValue $ws.Range('D51') '2.226.53'
Set-TextValue $ws.Range('E51') '  -1.29%  '
